# Applies the 2022-08-19 20:58:48 crawl refresh to the clothes sheet:
#  - insert the newly-observed "Selenacare Menstruationsunterwaesche S" row
#    (now on sale) at row 11, pushing the existing rows 11-13 down to 12-14
#  - drop the stale duplicate of that same product that used to sit at row 14
#    (after the insert it lives at row 15)
#  - flag the "Naturaline Herren T-Shirt Kurzarm schwarz M" listing as out of
#    stock online (row 4's aria-label)
#  - refresh every row's crawl timestamp (column O) to the new run time

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the new row for the on-sale Selenacare S product -------------
$ws.Rows("11").Insert()

# id, price and priceContextPrice are numeric-looking but must stay text
# (matches every other row in the sheet), so force a text format first.
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 11).NumberFormat = "@"

$ws.Cells.Item(11, 1).Value = "6548194"
$ws.Cells.Item(11, 2).Value = "Selenacare Menstruationsunterwäsche S"
$ws.Cells.Item(11, 3).Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/selenacare-menstruationsunterwaesche-s/p/6548194"
$ws.Cells.Item(11, 4).Value = "1ST"
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = "Selenacare"
$ws.Cells.Item(11, 8).Value = "12.25"
$ws.Cells.Item(11, 9).Value = "12.25/1ST"
$ws.Cells.Item(11, 10).Value = "Preis pro 1 Stück"
$ws.Cells.Item(11, 11).Value = "12.25"
$ws.Cells.Item(11, 12).Value = "1ST"
$ws.Cells.Item(11, 13).Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'unterwaesche']"
$ws.Cells.Item(11, 14).Value = "Selenacare Menstruationsunterwäsche S 50% Aktion 12.25 Schweizer Franken statt 24.50 Schweizer Franken"
$ws.Cells.Item(11, 15).Value = "2022-08-19 20:58:48"

# --- remove the now-duplicate old listing (shifted down to row 15) -------
$ws.Rows("15").Delete()

# --- out-of-stock annotation for the black M t-shirt ----------------------
$ws.Cells.Item(4, 14).Value = "Naturaline Herren T-Shirt Kurzarm schwarz M - Online kein Bestand 24.95 Schweizer Franken"

# --- refresh crawl timestamps on every data row (2-73) --------------------
$ws.Range("O2:O73").Value = "2022-08-19 20:58:48"
